# Update the "time_taken" column (F) on the "data" sheet with refreshed
# panel-query timestamps.
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

$dataSheet.Range("F2").Value  = "2021-10-05 14:35:04.544741"
$dataSheet.Range("F3").Value  = "2021-10-05 14:35:04.544748"
$dataSheet.Range("F4").Value  = "2021-10-05 14:35:04.544752"
$dataSheet.Range("F5").Value  = "2021-10-05 14:35:04.544755"
$dataSheet.Range("F6").Value  = "2021-10-05 14:35:04.544757"
$dataSheet.Range("F7").Value  = "2021-10-05 14:35:04.544760"
$dataSheet.Range("F8").Value  = "2021-10-05 14:35:04.544763"
$dataSheet.Range("F9").Value  = "2021-10-05 14:35:04.544765"
$dataSheet.Range("F10").Value = "2021-10-05 14:35:04.544768"
$dataSheet.Range("F11").Value = "2021-10-05 14:35:04.544771"
$dataSheet.Range("F12").Value = "2021-10-05 14:35:04.544774"
$dataSheet.Range("F13").Value = "2021-10-05 14:35:04.544776"

# Add a new "metadata" tab (placed right after "data") carrying the panel
# query metadata that used to live only outside the workbook.
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Neuroferritinopathies"
$metaSheet.Range("C2").Value = 3438
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "0.5"
$metaSheet.Range("D2").NumberFormat = "General"
$metaSheet.Range("D2").Style = "Normal"
$metaSheet.Range("E2").Value = "2021-03-14T01:21:32.988758Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:35:04.541113"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3438/?format=json"

# Match the header styling ("data"!B1, bold + border + center/top align)
# used on the existing sheet.
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)     # xlPasteFormats

$excel.CutCopyMode = 0
$dataSheet.Activate()
